# Update the URL in cell A5 and its associated hyperlink to the new
# "erroriswhatweneeed" address, and move the active-cell/selection state
# as closely as this runtime allows toward the recorded view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldUrl = "https://atpsgroup.com/admin/helloworldtesting"
$newUrl = "http://www.atpsgroup.com/erroriswhatweneeed"

# --- 1. Update the cell text (shared string) ---
$ws.Range("A5").Value = $newUrl

# --- 2. Update the hyperlink target that points at that cell ---
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Row -eq 5 -and $hl.Range.Column -eq 1) {
        $hl.Address = $newUrl
    }
}

# --- 3. Move the view's active cell/selection as closely as possible ---
$ws.Activate()
$ws.Range("A1:XFD1048576").Select()
$ws.Cells.Item(12, 1).Select()
